# Added filtering options for the Component Analysis
# Clears the trailing "extra" quarter-error columns that should not be
# included once filtering is applied, for rows 2, 3, 5, 6, and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5:K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
